$d = $word.ActiveDocument

# Apply each text replacement using a whole-document Find/Replace scoped to a single
# occurrence (wdReplaceOne). The replacements are ordered so that no intermediate
# document state contains an ambiguous / duplicate match for an "old" search string
# (this matters because one new value, 767×2=, equals another cell's original value,
# so the cell currently holding 767×2= must be updated to 723×7= before 812×8= is
# turned into a new 767×2=).

$r = $d.Content
$found = $r.Find.Execute("2024-08-24 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-25 Sunday", 1)
if (-not $found) {
    throw "Could not find expected text: 2024-08-24 Saturday"
}
$r = $d.Content
$found = $r.Find.Execute("403×4=", $true, $false, $false, $false, $false, $true, 1, $false, "518×6=", 1)
if (-not $found) {
    throw "Could not find expected text: 403×4="
}
$r = $d.Content
$found = $r.Find.Execute("767×2=", $true, $false, $false, $false, $false, $true, 1, $false, "723×7=", 1)
if (-not $found) {
    throw "Could not find expected text: 767×2="
}
$r = $d.Content
$found = $r.Find.Execute("812×8=", $true, $false, $false, $false, $false, $true, 1, $false, "767×2=", 1)
if (-not $found) {
    throw "Could not find expected text: 812×8="
}
$r = $d.Content
$found = $r.Find.Execute("146×9=", $true, $false, $false, $false, $false, $true, 1, $false, "931×6=", 1)
if (-not $found) {
    throw "Could not find expected text: 146×9="
}
$r = $d.Content
$found = $r.Find.Execute("620×5=", $true, $false, $false, $false, $false, $true, 1, $false, "632×7=", 1)
if (-not $found) {
    throw "Could not find expected text: 620×5="
}
$r = $d.Content
$found = $r.Find.Execute("403×7=", $true, $false, $false, $false, $false, $true, 1, $false, "624×3=", 1)
if (-not $found) {
    throw "Could not find expected text: 403×7="
}
$r = $d.Content
$found = $r.Find.Execute("623×4=", $true, $false, $false, $false, $false, $true, 1, $false, "218×3=", 1)
if (-not $found) {
    throw "Could not find expected text: 623×4="
}
$r = $d.Content
$found = $r.Find.Execute("217×6=", $true, $false, $false, $false, $false, $true, 1, $false, "702×9=", 1)
if (-not $found) {
    throw "Could not find expected text: 217×6="
}
$r = $d.Content
$found = $r.Find.Execute("303×8=", $true, $false, $false, $false, $false, $true, 1, $false, "826×7=", 1)
if (-not $found) {
    throw "Could not find expected text: 303×8="
}
$r = $d.Content
$found = $r.Find.Execute("912×8=", $true, $false, $false, $false, $false, $true, 1, $false, "278×4=", 1)
if (-not $found) {
    throw "Could not find expected text: 912×8="
}
$r = $d.Content
$found = $r.Find.Execute("836×3=", $true, $false, $false, $false, $false, $true, 1, $false, "476×5=", 1)
if (-not $found) {
    throw "Could not find expected text: 836×3="
}
$r = $d.Content
$found = $r.Find.Execute("105×4=", $true, $false, $false, $false, $false, $true, 1, $false, "932×8=", 1)
if (-not $found) {
    throw "Could not find expected text: 105×4="
}
$r = $d.Content
$found = $r.Find.Execute("440×2=", $true, $false, $false, $false, $false, $true, 1, $false, "209×5=", 1)
if (-not $found) {
    throw "Could not find expected text: 440×2="
}
$r = $d.Content
$found = $r.Find.Execute("754×9=", $true, $false, $false, $false, $false, $true, 1, $false, "461×7=", 1)
if (-not $found) {
    throw "Could not find expected text: 754×9="
}
$r = $d.Content
$found = $r.Find.Execute("185×2=", $true, $false, $false, $false, $false, $true, 1, $false, "549×4=", 1)
if (-not $found) {
    throw "Could not find expected text: 185×2="
}
$r = $d.Content
$found = $r.Find.Execute("841×5=", $true, $false, $false, $false, $false, $true, 1, $false, "829×3=", 1)
if (-not $found) {
    throw "Could not find expected text: 841×5="
}
$r = $d.Content
$found = $r.Find.Execute("206×9=", $true, $false, $false, $false, $false, $true, 1, $false, "495×5=", 1)
if (-not $found) {
    throw "Could not find expected text: 206×9="
}
$r = $d.Content
$found = $r.Find.Execute("574×5=", $true, $false, $false, $false, $false, $true, 1, $false, "482×5=", 1)
if (-not $found) {
    throw "Could not find expected text: 574×5="
}
$r = $d.Content
$found = $r.Find.Execute("780×4=", $true, $false, $false, $false, $false, $true, 1, $false, "139×6=", 1)
if (-not $found) {
    throw "Could not find expected text: 780×4="
}
$r = $d.Content
$found = $r.Find.Execute("580×3=", $true, $false, $false, $false, $false, $true, 1, $false, "231×7=", 1)
if (-not $found) {
    throw "Could not find expected text: 580×3="
}
$r = $d.Content
$found = $r.Find.Execute("306×8=", $true, $false, $false, $false, $false, $true, 1, $false, "322×4=", 1)
if (-not $found) {
    throw "Could not find expected text: 306×8="
}
$r = $d.Content
$found = $r.Find.Execute("965×3=", $true, $false, $false, $false, $false, $true, 1, $false, "329×7=", 1)
if (-not $found) {
    throw "Could not find expected text: 965×3="
}
$r = $d.Content
$found = $r.Find.Execute("976×3=", $true, $false, $false, $false, $false, $true, 1, $false, "916×4=", 1)
if (-not $found) {
    throw "Could not find expected text: 976×3="
}
$r = $d.Content
$found = $r.Find.Execute("854×3=", $true, $false, $false, $false, $false, $true, 1, $false, "465×5=", 1)
if (-not $found) {
    throw "Could not find expected text: 854×3="
}
$r = $d.Content
$found = $r.Find.Execute("992×2=", $true, $false, $false, $false, $false, $true, 1, $false, "446×7=", 1)
if (-not $found) {
    throw "Could not find expected text: 992×2="
}
